$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared-string table for the "Requisitos" detail rows (23-25) gets a new
# entry for LOM3229 inserted right after "Requisitos:", ahead of LOB1021 and
# LOM3016. Re-point the three rows so the rendered order becomes:
#   row 23 -> LOM3229 (Indicação de Conjunto)
#   row 24 -> LOB1021 (Requisito)
#   row 25 -> LOM3016 (Requisito)

$reqLOM3229 = "LOM3229 -  Métodos Experimentais da Física II  (Indicação de Conjunto)`n"
$reqLOB1021 = "LOB1021 -  Física IV  (Requisito)`n"
$reqLOM3016 = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"

$ws.Range("B23").Value = $reqLOM3229
$ws.Range("C23").Value = $reqLOM3229

$ws.Range("B24").Value = $reqLOB1021
$ws.Range("C24").Value = $reqLOB1021

$ws.Range("B25").Value = $reqLOM3016
$ws.Range("C25").Value = $reqLOM3016
